$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells remain text (matches source data which is
# stored as inline strings, not numbers) even for values that look numeric.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.783.64"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "1.770.18"
$ws.Range("E3").Value = "  -2.64%  "
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +1.42%  "
$ws.Range("D5").Value = "338.81"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("E6").Value = "  +1.13%  "
$ws.Range("D7").Value = "0.3798"
$ws.Range("E7").Value = "  -3.43%  "
$ws.Range("D8").Value = "0.3393"
$ws.Range("E8").Value = "  -2.78%  "
$ws.Range("D9").Value = "46.22"
$ws.Range("E9").Value = "  -4.88%  "
$ws.Range("D10").Value = "1.134"
$ws.Range("E10").Value = "  -6.30%  "
$ws.Range("D11").Value = "0.07340"
$ws.Range("E11").Value = "  -3.09%  "
$ws.Range("D12").Value = "23.25"
$ws.Range("E12").Value = "  +4.38%  "
$ws.Range("D13").Value = "1.006"
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("D14").Value = "6.313"
$ws.Range("E14").Value = "  -3.52%  "
$ws.Range("D15").Value = "7.365"
$ws.Range("E15").Value = "  +2.23%  "
$ws.Range("D16").Value = "1.773.76"
$ws.Range("E16").Value = "  -2.11%  "
$ws.Range("D17").Value = "0.00001067"
$ws.Range("E17").Value = "  -3.85%  "
$ws.Range("D18").Value = "0.06659"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("D19").Value = "81.78"
$ws.Range("E19").Value = "  -4.32%  "
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +1.07%  "
$ws.Range("D21").Value = "17.33"
$ws.Range("E21").Value = "  -2.98%  "
$ws.Range("D22").Value = "6.399"
$ws.Range("E22").Value = "  -2.88%  "
$ws.Range("D23").Value = "27.819.19"
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("D24").Value = "11.96"
$ws.Range("E24").Value = "  -7.21%  "
$ws.Range("D25").Value = "2.386"
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("D26").Value = "1.494"
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("D27").Value = "20.46"
$ws.Range("E27").Value = "  -4.74%  "
$ws.Range("D28").Value = "2.365"
$ws.Range("E28").Value = "  -8.03%  "
$ws.Range("D29").Value = "151.90"
$ws.Range("E29").Value = "  -2.17%  "
$ws.Range("D30").Value = "1.971.19"
$ws.Range("E30").Value = "  -2.11%  "
$ws.Range("D31").Value = "133.40"
$ws.Range("E31").Value = "  -1.69%  "
$ws.Range("D32").Value = "4.041"
$ws.Range("D33").Value = "5.961"
$ws.Range("E33").Value = "  -2.60%  "
$ws.Range("D34").Value = "0.08837"
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("D35").Value = "12.56"
$ws.Range("E35").Value = "  -5.79%  "
$ws.Range("D36").Value = "0.02379"
$ws.Range("E36").Value = "  -2.06%  "
$ws.Range("D37").Value = "0.6753"
$ws.Range("E37").Value = "  -2.66%  "
$ws.Range("D38").Value = "0.06328"
$ws.Range("E38").Value = "  -3.16%  "
$ws.Range("D39").Value = "5.234"
$ws.Range("E39").Value = "  -5.68%  "
$ws.Range("D40").Value = "0.2146"
$ws.Range("E40").Value = "  -3.74%  "
$ws.Range("D41").Value = "1.495"
$ws.Range("E41").Value = "  -7.42%  "
$ws.Range("D42").Value = "1.218"
$ws.Range("E42").Value = "  -3.66%  "
$ws.Range("D43").Value = "8.136"
$ws.Range("E43").Value = "  -5.02%  "
$ws.Range("E44").Value = "  +0.98%  "
$ws.Range("D45").Value = "14.02"
$ws.Range("E45").Value = "  -4.63%  "
$ws.Range("D46").Value = "0.6181"
$ws.Range("E46").Value = "  -5.76%  "
$ws.Range("D47").Value = "3.856"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").Value = "132.77"
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("D49").Value = "2.048"
$ws.Range("E49").Value = "  -5.51%  "
$ws.Range("D50").Value = "0.07341"
$ws.Range("E50").Value = "  +1.41%  "
$ws.Range("D51").Value = "1.197"
$ws.Range("E51").Value = "  +3.32%  "
